# "added pda screengrabs x2"
# Reposition / resize the body placeholder on slide 6 to make room for the
# newly added PDA screengrab images.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)

# EMU -> points conversion (1 pt = 12700 EMU)
$shp.Left = 361080 / 12700
$shp.Top  = 180000 / 12700
